# NATMI LR-pair results for Ccl21b-Ccr7 (YoungD7) were recomputed following
# Dr Hou's advice: the original single result row is replaced by four rows
# covering the Sending-cluster x Target-cluster combinations
# (FAPs/sCs x ECs/M2) now produced by the updated analysis.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> Ccl21b/Ccr7 -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ccl21b"
$ws.Range("C2").Value = "Ccr7"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2073196666666667
$ws.Range("H2").Value = 0.6219589999999999
$ws.Range("I2").Value = 0.4385567570045022
$ws.Range("J2").Value = 0.4385567570045022
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.6631983333333333
$ws.Range("N2").Value = 1.989595
$ws.Range("O2").Value = 0.4331413625787215
$ws.Range("P2").Value = 0.4331413625787215
$ws.Range("Q2").Value = 0.1374940574005555
$ws.Range("R2").Value = 1.237446516605
$ws.Range("S2").Value = 0.1899570712970353
$ws.Range("T2").Value = 0.1899570712970353

# Row 3: FAPs -> Ccl21b/Ccr7 -> M2
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ccl21b"
$ws.Range("C3").Value = "Ccr7"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2073196666666667
$ws.Range("H3").Value = 0.6219589999999999
$ws.Range("I3").Value = 0.4385567570045022
$ws.Range("J3").Value = 0.4385567570045022
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8679376666666666
$ws.Range("N3").Value = 2.603813
$ws.Range("O3").Value = 0.5668586374212784
$ws.Range("P3").Value = 0.5668586374212784
$ws.Range("Q3").Value = 0.1799405477407778
$ws.Range("R3").Value = 1.619464929667
$ws.Range("S3").Value = 0.2485996857074668
$ws.Range("T3").Value = 0.2485996857074668

# Row 4: sCs -> Ccl21b/Ccr7 -> ECs
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Ccl21b"
$ws.Range("C4").Value = "Ccr7"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.265412
$ws.Range("H4").Value = 0.796236
$ws.Range("I4").Value = 0.5614432429954979
$ws.Range("J4").Value = 0.5614432429954979
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.6631983333333333
$ws.Range("N4").Value = 1.989595
$ws.Range("O4").Value = 0.4331413625787215
$ws.Range("P4").Value = 0.4331413625787215
$ws.Range("Q4").Value = 0.1760207960466667
$ws.Range("R4").Value = 1.58418716442
$ws.Range("S4").Value = 0.2431842912816862
$ws.Range("T4").Value = 0.2431842912816862

# Row 5: sCs -> Ccl21b/Ccr7 -> M2
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Ccl21b"
$ws.Range("C5").Value = "Ccr7"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.265412
$ws.Range("H5").Value = 0.796236
$ws.Range("I5").Value = 0.5614432429954979
$ws.Range("J5").Value = 0.5614432429954979
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8679376666666666
$ws.Range("N5").Value = 2.603813
$ws.Range("O5").Value = 0.5668586374212784
$ws.Range("P5").Value = 0.5668586374212784
$ws.Range("Q5").Value = 0.2303610719853333
$ws.Range("R5").Value = 2.073249647868
$ws.Range("S5").Value = 0.3182589517138116
$ws.Range("T5").Value = 0.3182589517138116
